$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "confirm password"
$ws.Range("D1").Value = "message"

# Row 2
$ws.Range("A2").Value = "numpy84$3"
$ws.Range("B2").Value = "ninja"
$ws.Range("C2").Value = "ninja"
$ws.Range("D2").Value = "password_mismatch:The two password fields didn’t match."

# Row 3
$ws.Range("A3").Value = "numpy"
$ws.Range("B3").Value = "automation"
$ws.Range("C3").Value = "automation84"
$ws.Range("D3").Value = "password_mismatch:The two password fields didn’t match."

# Row 4
$ws.Range("A4").Value = "ninja"
$ws.Range("B4").Value = 123456789
$ws.Range("C4").Value = 123456789
$ws.Range("D4").Value = "password_mismatch:The two password fields didn’t match."

# Row 5
$ws.Range("A5").Value = "numpy"
$ws.Range("B5").Value = "ninja5"
$ws.Range("C5").Value = "ninja5"
$ws.Range("D5").Value = "password_mismatch:The two password fields didn’t match."

# Column width for the new "message" column
$ws.Columns.Item(4).ColumnWidth = 51.7

# Selection moves to C5
$ws.Range("C5").Select()
